$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "header" label that lived in C3, shifting the rest of
# row 3's content (E3 keeps "Empty rows ignored").
$ws.Range("C3").Clear()

# Move the active selection to C3 (the cell the author was inspecting).
$ws.Range("C3").Select()

# Set an explicit page setup, matching the document's print settings.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
